$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.946.83'
$ws.Range('E2').Value = '  +2.48%  '
$ws.Range('D3').Value = '2.422.51'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '552.13'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '137.86'
$ws.Range('E6').Value = '  +2.42%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.586'
$ws.Range('E8').Value = '  +2.71%  '
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.69'
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.148'
$ws.Range('E11').Value = '  -1.98%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.355'
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '25.45'
$ws.Range('E13').Value = '  +4.80%  '
$ws.Range('D14').Value = '2.856.11'
$ws.Range('E14').Value = '  +2.26%  '
$ws.Range('D15').Value = '59.855.91'
$ws.Range('E15').Value = '  +2.39%  '
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('D17').Value = '2.434.83'
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.36'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.41'
$ws.Range('E19').Value = '  +1.38%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '330.98'
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.71'
$ws.Range('E21').Value = '  -4.11%  '
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('E23').Value = '  +3.25%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.172'
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.77'
$ws.Range('E25').Value = '  +5.01%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  +3.52%  '
$ws.Range('D28').Value = '0.0₃0778'
$ws.Range('E28').Value = '  +4.04%  '
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '168.87'
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.14'
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '18.67'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('E33').Value = '  +1.48%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +3.44%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.22'
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '39.63'
$ws.Range('E39').Value = '  -1.79%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.412'
$ws.Range('E40').Value = '  -2.83%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '313.98'
$ws.Range('E41').Value = '  +6.76%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.68'
$ws.Range('E42').Value = '  -0.75%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '139.50'
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('E44').Value = '  +0.86%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0521'
$ws.Range('E45').Value = '  +0.83%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '19.50'
$ws.Range('E46').Value = '  +3.61%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.578'
$ws.Range('E47').Value = '  +2.09%  '
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.389'
$ws.Range('E49').Value = '  -5.33%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.66'
$ws.Range('E50').Value = '  +0.83%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '11.06'
$ws.Range('E51').Value = '  +0.21%  '
